$d = $word.ActiveDocument

# Update the date heading above the table.
$d.Content.Find.Execute("2025-07-26 Saturday", $true, $false, $false, $false, $false,
                         $true, 1, $false, "2025-07-27 Sunday", 2)

# Update the division problems in the table. Cells are addressed by
# (row, column) rather than by text search because several of the new
# values collide with old values elsewhere in the table (e.g. "63÷5="
# is both an original value and a value introduced by this edit), which
# would make a simple global find/replace ambiguous/unsafe.
$t = $d.Tables.Item(1)

$updates = @(
    @{Row = 1;  Col = 1; New = "46÷2="},
    @{Row = 1;  Col = 2; New = "55÷2="},
    @{Row = 1;  Col = 3; New = "29÷6="},
    @{Row = 1;  Col = 4; New = "81÷4="},
    @{Row = 1;  Col = 5; New = "76÷3="},

    @{Row = 5;  Col = 1; New = "21÷5="},
    @{Row = 5;  Col = 2; New = "12÷3="},
    @{Row = 5;  Col = 3; New = "96÷8="},
    @{Row = 5;  Col = 4; New = "60÷3="},
    @{Row = 5;  Col = 5; New = "43÷6="},

    @{Row = 9;  Col = 1; New = "42÷6="},
    @{Row = 9;  Col = 2; New = "89÷4="},
    @{Row = 9;  Col = 3; New = "74÷3="},
    @{Row = 9;  Col = 4; New = "89÷5="},
    @{Row = 9;  Col = 5; New = "11÷3="},

    @{Row = 13; Col = 1; New = "12÷7="},
    @{Row = 13; Col = 2; New = "66÷6="},
    @{Row = 13; Col = 3; New = "15÷9="},
    @{Row = 13; Col = 4; New = "63÷5="},
    @{Row = 13; Col = 5; New = "70÷4="},

    @{Row = 17; Col = 1; New = "44÷4="},
    @{Row = 17; Col = 2; New = "63÷5="},
    @{Row = 17; Col = 3; New = "47÷5="},
    @{Row = 17; Col = 4; New = "88÷9="},
    @{Row = 17; Col = 5; New = "48÷9="}
)

foreach ($u in $updates) {
    $cell = $t.Cell($u.Row, $u.Col)
    $cell.Range.Text = $u.New
}
